# Realestate Update resale numbers 2023-06-16 22:18
# Appends a new data row (row 51) to the CityResaleNum sheet with the
# latest resale-number snapshot, mirroring the existing rows' layout:
#   A: Date, B: Time, C: Weekday, D: Week (all stored as text)
#   E..T: numeric city values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51

# --- Text columns (A-D) -----------------------------------------------
# A plain Value assignment lets Excel auto-detect "date-looking" and
# "pure-integer-looking" strings and coerce them into numbers/dates.
# To keep these as literal text (matching the rest of the column), use
# a leading apostrophe to force text entry, then clear the resulting
# "quote prefix" cell format so the cell ends up looking like a normal,
# unstyled text cell (consistent with the other rows in the sheet).
$ws.Cells.Item($row, 1).Value = "'2023-06-16"
$ws.Cells.Item($row, 1).ClearFormats()

# Time-of-day text such as "22:12:53" and weekday names like "Friday"
# are not auto-converted by this engine, so they can be set directly.
$ws.Cells.Item($row, 2).Value = "22:12:53"
$ws.Cells.Item($row, 3).Value = "Friday"

$ws.Cells.Item($row, 4).Value = "'24"
$ws.Cells.Item($row, 4).ClearFormats()

# --- Numeric columns (E-T) ---------------------------------------------
$ws.Cells.Item($row, 5).Value  = 121859
$ws.Cells.Item($row, 6).Value  = 132999
$ws.Cells.Item($row, 7).Value  = 162097
$ws.Cells.Item($row, 8).Value  = 133143
$ws.Cells.Item($row, 9).Value  = 176996
$ws.Cells.Item($row, 10).Value = 114516
$ws.Cells.Item($row, 11).Value = 200884
$ws.Cells.Item($row, 12).Value = 224687
$ws.Cells.Item($row, 13).Value = 174890
$ws.Cells.Item($row, 14).Value = 103175
$ws.Cells.Item($row, 15).Value = 39189
$ws.Cells.Item($row, 16).Value = 34026
$ws.Cells.Item($row, 17).Value = 51735
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36889
$ws.Cells.Item($row, 20).Value = -1
